$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.00%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.82%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.555"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-12.81%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05890"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.76%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'-1.64%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8561"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.98%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9305"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-6.50%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01044"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1,626.81%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1406"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.48%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03589"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.16%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07083"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.54%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03232"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.48%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09213"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.38%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001555"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.05%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006099"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'5.16%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'0.58%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.195"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.54%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.63%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3061"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-3.49%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-0.98%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.852"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'8.94%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04209"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.10%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'0.18%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004292"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-6.17%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.04%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-22.02%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03837"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.75%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006237"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'15.18%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'-1.00%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-7.61%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01144"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'5.15%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005444"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.29%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.07197"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-15.75%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.1203"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'5,528.54%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("E50").Style = "Normal"
